$wb = $excel.ActiveWorkbook

# ALC row 41 (item id 5478)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1166.7333
$ws.Range("I41").Value = 1470.7
$ws.Range("J41").Value = 558.8
$ws.Range("K41").Value = 1470.7
$ws.Range("L41").Value = 558.8
$ws.Range("M41").Value = -1030.7
$ws.Range("N41").Value = -1438.8

# ALC row 64 (item id 5506)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3163.879
$ws.Range("I64").Value = 3257.0667
$ws.Range("J64").Value = 3086.2222
$ws.Range("K64").Value = 3257.0667
$ws.Range("L64").Value = 3086.2222
$ws.Range("M64").Value = -3009.0667
$ws.Range("N64").Value = -3582.2222

# ALC row 67 (item id 5506)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3163.879
$ws.Range("I67").Value = 3257.0667
$ws.Range("J67").Value = 3086.2222
$ws.Range("K67").Value = 3257.0667
$ws.Range("L67").Value = 3086.2222
$ws.Range("M67").Value = -2399.0667
$ws.Range("N67").Value = -4802.2222

# ALC row 70 (item id 12604)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1505.5111
$ws.Range("I70").Value = 1649.1666
$ws.Range("J70").Value = 1218.2
$ws.Range("K70").Value = 4947.4998
$ws.Range("L70").Value = 3654.6
$ws.Range("M70").Value = -4677.4998
$ws.Range("N70").Value = -4194.6

# ALC row 73 (item id 12604)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1505.5111
$ws.Range("I73").Value = 1649.1666
$ws.Range("J73").Value = 1218.2
$ws.Range("K73").Value = 4947.4998
$ws.Range("L73").Value = 3654.6
$ws.Range("M73").Value = -4011.4998
$ws.Range("N73").Value = -5526.6

# ALC row 135 (item id 44047)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1083.129
$ws.Range("I135").Value = 775.5925999999999
$ws.Range("J135").Value = 3159
$ws.Range("K135").Value = 6980.3334
$ws.Range("L135").Value = 28431
$ws.Range("M135").Value = -4445.3334
$ws.Range("N135").Value = -33501

# ALC row 138 (item id 44169)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3064.7942
$ws.Range("I138").Value = 1891.8
$ws.Range("J138").Value = 3396.7737
$ws.Range("K138").Value = 5675.4
$ws.Range("L138").Value = 10190.3211
$ws.Range("M138").Value = -535.3999999999996
$ws.Range("N138").Value = -20470.3211

# ARM row 2 (item id 27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 943.4865
$ws.Range("I2").Value = 824
$ws.Range("J2").Value = 1100.3125
$ws.Range("K2").Value = 824
$ws.Range("L2").Value = 1100.3125
$ws.Range("M2").Value = -711
$ws.Range("N2").Value = -1326.3125

# ARM row 61 (item id 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4015
$ws.Range("I61").Value = 2618.2
$ws.Range("J61").Value = 7507
$ws.Range("K61").Value = 2618.2
$ws.Range("L61").Value = 7507
$ws.Range("M61").Value = -2406.2
$ws.Range("N61").Value = -7931

# ARM row 63 (item id 12528)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3648.5
$ws.Range("I63").Value = 2137.8
$ws.Range("J63").Value = 5159.2
$ws.Range("K63").Value = 2137.8
$ws.Range("L63").Value = 5159.2
$ws.Range("M63").Value = -1451.8
$ws.Range("N63").Value = -6531.2

# ARM row 66 (item id 12528)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3648.5
$ws.Range("I66").Value = 2137.8
$ws.Range("J66").Value = 5159.2
$ws.Range("K66").Value = 10689
$ws.Range("L66").Value = 25796
$ws.Range("M66").Value = -7257
$ws.Range("N66").Value = -32660

# ARM row 116 (item id 27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 943.4865
$ws.Range("I116").Value = 824
$ws.Range("J116").Value = 1100.3125
$ws.Range("K116").Value = 824
$ws.Range("L116").Value = 1100.3125
$ws.Range("M116").Value = 1470
$ws.Range("N116").Value = -5688.3125

# ARM row 136 (item id 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4015
$ws.Range("I136").Value = 2618.2
$ws.Range("J136").Value = 7507
$ws.Range("K136").Value = 7854.599999999999
$ws.Range("L136").Value = 22521
$ws.Range("M136").Value = -5304.599999999999
$ws.Range("N136").Value = -27621

# BSM row 3 (item id 27713)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 943.4865
$ws.Range("I3").Value = 824
$ws.Range("J3").Value = 1100.3125
$ws.Range("K3").Value = 824
$ws.Range("L3").Value = 1100.3125
$ws.Range("M3").Value = -710
$ws.Range("N3").Value = -1328.3125

# BSM row 22 (item id 5092)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 125
$ws.Range("I22").Value = 125
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 125
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 48
$ws.Range("N22").ClearContents()

# BSM row 94 (item id 19939)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 7626.5
$ws.Range("I94").Value = 828.15
$ws.Range("J94").Value = 21223.2
$ws.Range("K94").Value = 828.15
$ws.Range("L94").Value = 21223.2
$ws.Range("M94").Value = -377.15
$ws.Range("N94").Value = -22125.2

# BSM row 107 (item id 27706)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1431.75
$ws.Range("I107").Value = 1353.9286
$ws.Range("J107").Value = 1613.3334
$ws.Range("K107").Value = 1353.9286
$ws.Range("L107").Value = 1613.3334
$ws.Range("M107").Value = 566.0714
$ws.Range("N107").Value = -5453.3334

# BSM row 134 (item id 43998)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2050.842
$ws.Range("I134").Value = 1463.8334
$ws.Range("J134").Value = 3057.1428
$ws.Range("K134").Value = 4391.5002
$ws.Range("L134").Value = 9171.428400000001
$ws.Range("M134").Value = -1856.5002
$ws.Range("N134").Value = -14241.4284

# CRP row 5 (item id 1893)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 206.5
$ws.Range("I5").Value = 109.28571
$ws.Range("J5").Value = 433.33334
$ws.Range("K5").Value = 109.28571
$ws.Range("L5").Value = 433.33334
$ws.Range("M5").Value = 2.714290000000005
$ws.Range("N5").Value = -657.33334

# CRP row 29 (item id 2408)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 22021
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 22021
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 22021
$ws.Range("N29").Value = -22607

# CRP row 31 (item id 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7845464.5
$ws.Range("I31").Value = 18182924
$ws.Range("J31").Value = 3254.4827
$ws.Range("K31").Value = 18182924
$ws.Range("L31").Value = 3254.4827
$ws.Range("M31").Value = -18182629
$ws.Range("N31").Value = -3844.4827

# CRP row 34 (item id 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 7845464.5
$ws.Range("I34").Value = 18182924
$ws.Range("J34").Value = 3254.4827
$ws.Range("K34").Value = 18182924
$ws.Range("L34").Value = 3254.4827
$ws.Range("M34").Value = -18182722
$ws.Range("N34").Value = -3658.4827

# CUL row 4 (item id 4650)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 159.35294
$ws.Range("I4").Value = 159.35294
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 478.05882
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -366.05882

# CUL row 5 (item id 43974)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 754.1163
$ws.Range("I5").Value = 498.72
$ws.Range("J5").Value = 1108.8334
$ws.Range("K5").Value = 1496.16
$ws.Range("L5").Value = 3326.5002
$ws.Range("M5").Value = -1384.16
$ws.Range("N5").Value = -3550.5002

# CUL row 113 (item id 27843)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 984.625
$ws.Range("I113").Value = 422.69446
$ws.Range("J113").Value = 1707.1072
$ws.Range("K113").Value = 1268.08338
$ws.Range("L113").Value = 5121.321599999999
$ws.Range("M113").Value = 901.91662
$ws.Range("N113").Value = -9461.321599999999

# CUL row 135 (item id 43974)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 754.1163
$ws.Range("I135").Value = 498.72
$ws.Range("J135").Value = 1108.8334
$ws.Range("K135").Value = 4488.48
$ws.Range("L135").Value = 9979.500599999999
$ws.Range("M135").Value = -1953.48
$ws.Range("N135").Value = -15049.5006

# GSM row 14 (item id 4198)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 4168.3335
$ws.Range("I14").Value = 500
$ws.Range("J14").Value = 6002.5
$ws.Range("K14").Value = 500
$ws.Range("L14").Value = 6002.5
$ws.Range("M14").Value = -332
$ws.Range("N14").Value = -6338.5

# GSM row 20 (item id 4095)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()

# GSM row 113 (item id 27710)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 916.3226
$ws.Range("I113").Value = 788.52
$ws.Range("J113").Value = 1448.8334
$ws.Range("K113").Value = 788.52
$ws.Range("L113").Value = 1448.8334
$ws.Range("M113").Value = 1381.48
$ws.Range("N113").Value = -5788.8334

# GSM row 135 (item id 42006)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 18000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 18000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 18000
$ws.Range("N135").Value = -28140

# GSM row 138 (item id 42325)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 41533.332
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 41533.332
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 41533.332
$ws.Range("N138").Value = -51813.332

# LTW row 44 (item id 3658)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H44").Value = 12400
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 12400
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 12400
$ws.Range("N44").Value = -13312

# LTW row 137 (item id 43296)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H137").Value = 19716.666
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 19716.666
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 19716.666
$ws.Range("N137").Value = -29916.666

# WVR row 59 (item id 3201)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

# WVR row 113 (item id 27752)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 563.8333
$ws.Range("I113").Value = 491.625
$ws.Range("J113").Value = 708.25
$ws.Range("K113").Value = 1474.875
$ws.Range("L113").Value = 2124.75
$ws.Range("M113").Value = 695.125
$ws.Range("N113").Value = -6464.75
